$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 24 (Excel's native Insert semantics
# shift rows 24:39 down to 26:41, preserving formatting/styles of row 24).
$ws.Rows("24:25").Insert()

# ---- New row 24 ----
$ws.Cells.Item(24, 1).Value = 9
$ws.Cells.Item(24, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44421
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100102
$ws.Cells.Item(24, 8).Value = "Cítricos"
$ws.Cells.Item(24, 9).Value = 100102006
$ws.Cells.Item(24, 10).Value = "Pomelo"
$ws.Cells.Item(24, 11).Value = "Start Ruby"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 80
$ws.Cells.Item(24, 14).Value = 8400
$ws.Cells.Item(24, 15).Value = 8400
$ws.Cells.Item(24, 16).Value = 8400
$ws.Cells.Item(24, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(24, 18).Value = "Región Metropolitana"
$ws.Cells.Item(24, 19).Value = 600
$ws.Cells.Item(24, 20).Value = 14

# ---- New row 25 ----
$ws.Cells.Item(25, 1).Value = 9
$ws.Cells.Item(25, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44421
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100102
$ws.Cells.Item(25, 8).Value = "Cítricos"
$ws.Cells.Item(25, 9).Value = 100102006
$ws.Cells.Item(25, 10).Value = "Pomelo"
$ws.Cells.Item(25, 11).Value = "Start Ruby"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 120
$ws.Cells.Item(25, 14).Value = 7000
$ws.Cells.Item(25, 15).Value = 7000
$ws.Cells.Item(25, 16).Value = 7000
$ws.Cells.Item(25, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(25, 18).Value = "Región Metropolitana"
$ws.Cells.Item(25, 19).Value = 500
$ws.Cells.Item(25, 20).Value = 14
